$d = $word.ActiveDocument

# 1) Italic English "Programa resumido" sentence: insert a line break after
#    "Equilibrium of Rigid Bodies, "
$d.Content.Find.Execute(
    "Statics of Particles, Statics of Rigid Bodies, Equilibrium of Rigid Bodies, Analysis of Structures, Distributed Forces, Forces in Beams.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Statics of Particles, Statics of Rigid Bodies, Equilibrium of Rigid Bodies, ^lAnalysis of Structures, Distributed Forces, Forces in Beams.",
    2)

# 2) Portuguese "Programa" paragraph: break into numbered items
$d.Content.Find.Execute(
    "Programa: 1) Estática de partículas: Vetores. Resultante de várias forças concorrentes. Equilíbrio de uma partícula. 2) Estática de Corpos Rígidos: Conceito de corpo rígido. Forças externas e forças internas. Forças equivalentes. Momento de uma força com relação a um ponto. Sistemas equivalentes de forças. Diagrama de corpo livre.3) Equilíbrio de corpos rígidos: Reações de apoios e conexões para uma estrutura 2D. Equilíbrio de um corpo rígido em 2D. Reações de apoios e conexões para uma estrutura 3D. Equilíbrio de um corpo rígido em 3D. 4) Análise de Estruturas: Treliças: Definições. Treliça simples. Análise de treliças pelo método dos nós. Análise de treliças pelo método das seções. Estruturas: estruturas que contêm elementos sujeitos a ação de múltiplas forças, transmissão e modificação de forças.5) Forças Distribuídas: Centróides e baricentros de áreas, linhas e volumes; teoremas de Guldinus-Pappus; cargas distribuídas sobre vigas.6) Forças em Vigas: Diagramas de cisalhamento e momento fletor para uma viga carregada.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Programa: ^l^l1) Estática de partículas: Vetores. Resultante de várias forças concorrentes. Equilíbrio de uma partícula. ^l2) Estática de Corpos Rígidos: Conceito de corpo rígido. Forças externas e forças internas. Forças equivalentes. Momento de uma força com relação a um ponto. Sistemas equivalentes de forças. Diagrama de corpo livre.^l3) Equilíbrio de corpos rígidos: Reações de apoios e conexões para uma estrutura 2D. Equilíbrio de um corpo rígido em 2D. Reações de apoios e conexões para uma estrutura 3D. Equilíbrio de um corpo rígido em 3D. ^l4) Análise de Estruturas: Treliças: Definições. Treliça simples. Análise de treliças pelo método dos nós. Análise de treliças pelo método das seções. Estruturas: estruturas que contêm elementos sujeitos a ação de múltiplas forças, transmissão e modificação de forças.^l5) Forças Distribuídas: Centróides e baricentros de áreas, linhas e volumes; teoremas de Guldinus-Pappus; cargas distribuídas sobre vigas.^l^l6) Forças em Vigas: Diagramas de cisalhamento e momento fletor para uma viga carregada.",
    2)

# 3) Italic English "Programa" paragraph: break into numbered items
$d.Content.Find.Execute(
    "1) Static particles: Vectors. Resulting from various competing forces. Equilibrium of a particle.2) Statics of Rigid Bodies: Rigid body concept. External forces and internal forces. Equivalent forces. Moment of a force with respect to a point. Systems equivalent forces. Free body diagram.3) Equilibrium of rigid bodies: Support reactions and connections to a 2D structure. Equilibrium of a rigid body 2D. Support reactions and connections to a 3D structure. Equilibrium of a rigid body in 3D.4) Analysis of Structures: Trusses : Definitions . Simple trusses . Trusses analysis by the method of nodes. Trusses analysis by the method of sections. Structures: structures that contain elements subject to action of multiple forces , transmission and modification forces.5) Distributed Forces: Barycentres, centroids and areas, lines and volumes; theorems of Pappus-Guldinus, distributed loads on beams.6) Forces in Beams: Diagrams of shear and bending moment for a beam loaded.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1) Static particles: Vectors. Resulting from various competing forces. Equilibrium of a particle.^l2) Statics of Rigid Bodies: Rigid body concept. External forces and internal forces. Equivalent forces. Moment of a force with respect to a point. Systems equivalent forces. Free body diagram.^l3) Equilibrium of rigid bodies: Support reactions and connections to a 2D structure. Equilibrium of a rigid body 2D. Support reactions and connections to a 3D structure. Equilibrium of a rigid body in 3D.^l4) Analysis of Structures: Trusses : Definitions . Simple trusses . Trusses analysis by the method of nodes. Trusses analysis by the method of sections. Structures: structures that contain elements subject to action of multiple forces , transmission and modification forces.^l5) Distributed Forces: Barycentres, centroids and areas, lines and volumes; theorems of Pappus-Guldinus, distributed loads on beams.^l6) Forces in Beams: Diagrams of shear and bending moment for a beam loaded.",
    2)

# 4) Bibliography paragraph: break into numbered references
$d.Content.Find.Execute(
    "1. BEER, Ferdinand Pierre, ; JOHNSTON, E. Russel.; Eisenberg, Elliot R., Mecânica vetorial para engenheiros: Estática.  Mc Graw Hill (2011).2. HIBBELER, R.C. Mecânica para engenharia, Vol. 1: estática, Pearson Prentice Hall (2005).3. MERIAM J. L. ; KRAIGE, L. G., Mecânica, estática, Livros Técnicos e Científicos Editora (2004).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. BEER, Ferdinand Pierre, ; JOHNSTON, E. Russel.; Eisenberg, Elliot R., Mecânica vetorial para engenheiros: Estática.  Mc Graw Hill (2011).^l2. HIBBELER, R.C. Mecânica para engenharia, Vol. 1: estática, Pearson Prentice Hall (2005).^l3. MERIAM J. L. ; KRAIGE, L. G., Mecânica, estática, Livros Técnicos e Científicos Editora (2004).",
    2)
